# Lotofacil simulator - append the latest contest draws (rows 302-307)
# and adjust the sheet view, matching the author's layout tweak.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Lotofacil results to append below the existing data (A1:P301 -> A1:P307)
$novosSorteios = @(
    @(302, 3466, 1, 5, 6, 7, 9, 10, 12, 15, 17, 18, 19, 21, 23, 24, 25),
    @(303, 3467, 3, 4, 5, 6, 7, 8, 9, 13, 14, 16, 18, 19, 20, 21, 25),
    @(304, 3468, 1, 3, 4, 6, 9, 10, 12, 13, 14, 18, 19, 21, 23, 24, 25),
    @(305, 3469, 1, 3, 5, 6, 7, 8, 11, 12, 13, 15, 17, 18, 21, 22, 24),
    @(306, 3470, 1, 4, 5, 7, 8, 10, 12, 13, 14, 18, 20, 21, 22, 23, 24),
    @(307, 3471, 1, 2, 3, 5, 6, 11, 13, 16, 17, 19, 21, 22, 23, 24, 25)
)

foreach ($sorteio in $novosSorteios) {
    $linha = $sorteio[0]
    for ($col = 1; $col -le 16; $col++) {
        $ws.Cells.Item($linha, $col).Value = $sorteio[$col]
    }
}

# Scroll the view down toward the newly added rows and select the last row,
# mirroring where the author left the cursor after pasting the new data.
$win = $excel.ActiveWindow
$win.ScrollRow = 287
$win.ScrollColumn = 1
$ws.Rows.Item(301).Select()
